# Update "想去人数" (F column) and "最低票价" (G column) figures across sheets
# to reflect newly generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 593
$ws.Range("F4").Value = 6392
$ws.Range("G4").Value = 70
$ws.Range("F5").Value = 723
$ws.Range("F8").Value = 324
$ws.Range("F11").Value = 697
$ws.Range("F12").Value = 1179
$ws.Range("F14").Value = 427
$ws.Range("F17").Value = 1418
$ws.Range("F20").Value = 399
$ws.Range("F23").Value = 147
$ws.Range("F24").Value = 2216
$ws.Range("F26").Value = 97
$ws.Range("F29").Value = 3576
$ws.Range("F31").Value = 633

# --- Sheet: 演出 (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 181
$ws.Range("F9").Value = 108
$ws.Range("F10").Value = 20
$ws.Range("F11").Value = 1018
$ws.Range("F13").Value = 105
$ws.Range("F16").Value = 75
$ws.Range("F21").Value = 13
$ws.Range("F22").Value = 15
$ws.Range("F25").Value = 231
$ws.Range("F26").Value = 88
$ws.Range("F34").Value = 3

# --- Sheet: 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1194
$ws.Range("F6").Value = 1575
$ws.Range("F10").Value = 771

# --- Sheet: 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1194
$ws.Range("F4").Value = 1575
$ws.Range("F7").Value = 771
$ws.Range("F8").Value = 593
$ws.Range("F10").Value = 6392
$ws.Range("G10").Value = 70
$ws.Range("F12").Value = 723
$ws.Range("F16").Value = 324
$ws.Range("F18").Value = 697
$ws.Range("F19").Value = 105
$ws.Range("F20").Value = 75
$ws.Range("F21").Value = 1179
$ws.Range("F23").Value = 427
$ws.Range("F27").Value = 1418
$ws.Range("F28").Value = 13
$ws.Range("F31").Value = 399
$ws.Range("F34").Value = 231
$ws.Range("F36").Value = 147
$ws.Range("F38").Value = 2216
$ws.Range("F42").Value = 97
$ws.Range("F45").Value = 3576
$ws.Range("F46").Value = 3
$ws.Range("F49").Value = 633
